$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format only on Price cells whose new value would otherwise be
# auto-parsed by Excel as a number, so they stay text like the original cells.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.265.81'
$ws.Range('E2').Value = '  -4.65%  '
$ws.Range('D3').Value = '1.752.53'
$ws.Range('E3').Value = '  -4.40%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '1.003'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '306.48'
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('D7').Value = '0.4275'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.3603'
$ws.Range('E8').Value = '  -1.71%  '
$ws.Range('D9').Value = '0.07048'
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('D10').Value = '0.8252'
$ws.Range('E10').Value = '  -4.38%  '
$ws.Range('D11').Value = '19.94'
$ws.Range('E11').Value = '  -3.27%  '
$ws.Range('D12').Value = '1.754.18'
$ws.Range('E12').Value = '  -4.53%  '
$ws.Range('D13').Value = '5.176'
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').Value = '6.290'
$ws.Range('E14').Value = '  -3.43%  '
$ws.Range('D15').Value = '0.06824'
$ws.Range('E15').Value = '  -1.70%  '
$ws.Range('D16').Value = '1.009'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').Value = '78.41'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').Value = '0.000008607'
$ws.Range('E18').Value = '  -3.31%  '
$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').Value = '14.77'
$ws.Range('E20').Value = '  -3.86%  '
$ws.Range('D21').Value = '26.481.35'
$ws.Range('E21').Value = '  -4.00%  '
$ws.Range('D22').Value = '4.936'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('D23').Value = '11.05'
$ws.Range('E23').Value = '  +2.12%  '
$ws.Range('D24').Value = '2.000.04'
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('D25').Value = '1.897'
$ws.Range('E25').Value = '  -4.79%  '
$ws.Range('D26').Value = '151.77'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').Value = '17.98'
$ws.Range('E27').Value = '  -4.59%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').Value = '114.21'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '4.981'
$ws.Range('E29').Value = '  -2.45%  '
$ws.Range('D30').Value = '1.619'
$ws.Range('E30').Value = '  -11.04%  '
$ws.Range('D31').Value = '0.08891'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').Value = '0.7068'
$ws.Range('E32').Value = '  -5.33%  '
$ws.Range('D33').Value = '4.259'
$ws.Range('E33').Value = '  -6.15%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.081'
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = '1.004'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.749'
$ws.Range('E36').Value = '  -7.97%  '
$ws.Range('D37').Value = '1.060'
$ws.Range('E37').Value = '  -2.72%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.05052'
$ws.Range('E38').Value = '  -4.75%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01868'
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('D40').Value = '0.4847'
$ws.Range('E40').Value = '  -4.43%  '
$ws.Range('D41').Value = '0.1585'
$ws.Range('E41').Value = '  -4.54%  '
$ws.Range('D42').Value = '2.474'
$ws.Range('E42').Value = '  -11.55%  '
$ws.Range('D43').Value = '6.110'
$ws.Range('E43').Value = '  -5.99%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '7.827'
$ws.Range('E45').Value = '  -5.72%  '
$ws.Range('D46').Value = '103.95'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = '9.956'
$ws.Range('E47').Value = '  -4.12%  '
$ws.Range('D48').Value = '0.06171'
$ws.Range('E48').Value = '  -4.79%  '
$ws.Range('D49').Value = '0.4415'
$ws.Range('E49').Value = '  -5.48%  '
$ws.Range('D50').Value = '1.551'
$ws.Range('E50').Value = '  -3.65%  '
$ws.Range('D51').Value = '1.693'
